$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("B2").Value = 6.26
$ws.Range("A3").Value = -21.663
$ws.Range("C3").Value = -12.774
$ws.Range("C12").Value = -11.288
$ws.Range("A14").Value = -21.603
$ws.Range("A21").Value = -20.257
$ws.Range("A23").Value = -20.585
$ws.Range("C24").Value = -12.679
$ws.Range("A25").Value = -20.53900000000001
$ws.Range("B25").Value = 7.273999999999999
$ws.Range("C25").Value = -12.961
$ws.Range("A26").Value = -21.326
$ws.Range("B27").Value = 5.513
$ws.Range("A29").Value = -21.219
$ws.Range("B31").Value = 5.988
$ws.Range("B39").Value = 7.581999999999999
$ws.Range("B48").Value = 5.24
$ws.Range("C50").Value = -13.133
$ws.Range("B51").Value = 5.42
$ws.Range("B52").Value = 5.44
$ws.Range("A53").Value = -21.841
$ws.Range("C53").Value = -11.551
$ws.Range("B55").Value = 4.583
$ws.Range("B56").Value = 4.907
$ws.Range("A57").Value = -21.352
$ws.Range("B57").Value = 5.948
$ws.Range("C57").Value = -13.257
$ws.Range("A59").Value = -22.5
$ws.Range("C61").Value = -13.025
$ws.Range("C63").Value = -11.518
$ws.Range("A69").Value = -21.507
$ws.Range("C70").Value = -12.045
$ws.Range("B73").Value = 6.439
$ws.Range("A79").Value = -21.246
$ws.Range("A83").Value = -22.015
$ws.Range("C86").Value = -13.252
$ws.Range("B89").Value = 5.915999999999999
$ws.Range("B90").Value = 5.907
$ws.Range("A91").Value = -21.533
$ws.Range("B92").Value = 6.027
$ws.Range("A93").Value = -21.439
$ws.Range("C98").Value = -12.45
$ws.Range("C100").Value = -13.247
$ws.Range("C102").Value = -13.361
